$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
# This shared string is referenced from the Overview sheet (zh-cn/de-de
# status columns E & F) as well as from each language sheet's Status
# column (C). Every cell that currently holds the old text needs to be
# updated so they keep sharing the (renamed) string.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the "Status" columns on every sheet ---
# Overview: columns E (zh-cn) and F (de-de)
$wsOverview.Columns.Item(5).ColumnWidth = 12.45
$wsOverview.Columns.Item(6).ColumnWidth = 12.45

# zh-cn / de-de: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.45
$wsDeDe.Columns.Item(3).ColumnWidth = 12.45
